$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.512.56'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '1.877.28'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("D4").Value = '1.025'
$ws.Range("E4").Value = '  +1.99%  '
$ws.Range("D5").Value = '317.49'
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").Value = '1.022'
$ws.Range("E6").Value = '  +1.68%  '
$ws.Range("D7").Value = '0.5148'
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").Value = '0.3941'
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("D9").Value = '0.08327'
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("D10").Value = '1.118'
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").Value = '42.10'
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("D12").Value = '6.248'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '1.870.99'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("D14").Value = '20.47'
$ws.Range("E14").Value = '  -1.10%  '
$ws.Range("D15").Value = '7.246'
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").Value = '1.025'
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").Value = '0.00001112'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '91.35'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '0.06761'
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("D20").Value = '17.74'
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").Value = '1.022'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").Value = '5.975'
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").Value = '28.541.17'
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("D24").Value = '11.20'
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").Value = '2.268'
$ws.Range("E25").Value = '  +0.08%  '
$ws.Range("D26").Value = '2.079.16'
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").Value = '161.75'
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("D28").Value = '20.71'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").Value = '2.380'
$ws.Range("E29").Value = '  -4.91%  '
$ws.Range("D30").Value = '127.54'
$ws.Range("E30").Value = '  +1.95%  '
$ws.Range("D31").Value = '0.1055'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("D32").Value = '1.038'
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").Value = '5.842'
$ws.Range("E33").Value = '  -0.60%  '
$ws.Range("D34").Value = '3.660'
$ws.Range("E34").Value = '  +1.68%  '
$ws.Range("D35").Value = '0.02441'
$ws.Range("E35").Value = '  -0.81%  '
$ws.Range("D36").Value = '0.06513'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").Value = '9.148'
$ws.Range("E37").Value = '  -5.68%  '
$ws.Range("D38").Value = '0.2190'
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").Value = '1.253'
$ws.Range("E39").Value = '  +1.60%  '
$ws.Range("D40").Value = '1.192'
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").Value = '0.6451'
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("D42").Value = '5.008'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '11.20'
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("D44").Value = '0.6038'
$ws.Range("E44").Value = '  -1.72%  '
$ws.Range("D45").Value = '12.95'
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").Value = '3.719'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("D47").Value = '1.251'
$ws.Range("E47").Value = '  -2.76%  '
$ws.Range("D48").Value = '2.001'
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").Value = '1.215'
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").Value = '122.16'
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("D51").Value = '0.06876'
$ws.Range("E51").Value = '  -0.34%  '
